# Auto-generated edit script for cryptos.xlsx price/volume update
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '70.439.56'
$ws.Range('E2').Value = '  +0.83%  '
$ws.Range('D3').Value = '3.625.80'
$ws.Range('E3').Value = '  +2.70%  '
$ws.Range('E4').Value = '  +0.16%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '601.95'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.53%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '196.00'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.56%  '
$ws.Range('E7').Value = '  -0.44%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.999'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +0.06%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.211'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +5.01%  '
$ws.Range('E10').Value = '  -0.14%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '53.35'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.67%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0000303'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -0.13%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '9.59'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +0.72%  '
$ws.Range('D14').Value = '4.200.44'
$ws.Range('E14').Value = '  +2.78%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '606.36'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +1.52%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '12.95'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +1.81%  '
$ws.Range('D17').Value = '70.537.01'
$ws.Range('E17').Value = '  +0.80%  '
$ws.Range('B18').Value = 'Chainlink'
$ws.Range('C18').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '19.12'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +0.07%  '
$ws.Range('B19').Value = 'WrappedEther'
$ws.Range('C19').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D19').Value = '3.617.15'
$ws.Range('E19').Value = '  +2.84%  '
$ws.Range('E20').Value = '  +1.42%  '
$ws.Range('E21').Value = '  +0.67%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '18.25'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -1.04%  '
$ws.Range('E23').Value = '  -1.38%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '103.44'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +1.37%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '4.62'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -0.27%  '
$ws.Range('E26').Value = '  -5.48%  '
$ws.Range('E27').Value = '  -2.13%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '9.70'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +0.99%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '33.92'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +1.23%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '4.67'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +8.22%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '7.33'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +3.22%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '12.29'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -1.73%  '
$ws.Range('E33').Value = '  +1.86%  '
$ws.Range('E34').Value = '  +0.51%  '
$ws.Range('D35').Value = '0.0₃0883'
$ws.Range('E35').Value = '  +2.78%  '
$ws.Range('D36').Value = '3.910.54'
$ws.Range('E36').Value = '  +4.79%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.00'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +0.03%  '
$ws.Range('B38').Value = 'Fetch.AI'
$ws.Range('C38').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '3.07'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -0.12%  '
$ws.Range('B39').Value = 'Bittensor'
$ws.Range('C39').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '520.19'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +6.20%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '36.97'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +0.98%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.388'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -1.05%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '3.55'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -2.39%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.135'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +1.02%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0460'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +1.01%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '3.61'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +8.94%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.92'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +2.94%  '
$ws.Range('E47').Value = '  -0.28%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '8.57'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +0.29%  '
$ws.Range('E49').Value = '  -0.18%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.000249'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -1.64%  '
$ws.Range('E51').Value = '  +1.50%  '
